$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new timesheet entry row (row 10): date + hours worked, reusing the
# date formatting already applied to column A (rows 1-9) by copying A9's
# formatting into A10 before writing the new values.
$ws.Range("A9").Copy($ws.Range("A10"))
$ws.Range("A10").Value = 44350
$ws.Range("B10").Value = 1

# Extend the selection to cover the new row, same as the original sheet's
# B1:B9 selection growing to B1:B10.
$ws.Range("B1").Select()
$ws.Range("B1:B10").Select()
